$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.749.84'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '2.678.80'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'601.11"
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').Value = "'156.90"
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +5.82%  '
$ws.Range('E9').Value = '  +5.61%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').Value = "'5.89"
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').Value = "'0.401"
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = "'29.44"
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '3.159.69'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '65.587.24'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '2.670.05'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').Value = "'12.96"
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').Value = "'7.61"
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = "'352.81"
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = "'69.77"
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  +5.95%  '
$ws.Range('D25').Value = "'9.70"
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('E28').Value = '  -5.88%  '
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = "'530.95"
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('D33').Value = "'1.75"
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').Value = "'6.50"
$ws.Range('E34').Value = '  -3.37%  '
$ws.Range('D35').Value = "'5.50"
$ws.Range('E35').Value = '  +0.63%  '
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').Value = "'20.59"
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = "'158.42"
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('E40').Value = '  -2.87%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'165.02"
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').Value = "'22.93"
$ws.Range('E46').Value = '  -2.82%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = "'0.0259"
$ws.Range('E47').Value = '  -2.96%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.642"
$ws.Range('E48').Value = '  -2.64%  '
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('D51').Value = "'20.13"
$ws.Range('E51').Value = '  -4.37%  '
